$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$strs = @{20="ECs"; 21="FAPs"; 22="sCs"; 23="Pros1"; 24="Tyro3"}

$data = New-Object 'object[,]' 9,20
$data[0,0] = $strs[20]
$data[0,1] = $strs[23]
$data[0,2] = $strs[24]
$data[0,3] = $strs[20]
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 17.586354
$data[0,7] = 52.759062
$data[0,8] = 0.2178245326054132
$data[0,9] = 0.2178245326054132
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.07610666666666667
$data[0,13] = 0.22832
$data[0,14] = 0.04153709822108353
$data[0,15] = 0.04153709822108354
$data[0,16] = 1.33843878176
$data[0,17] = 12.04594903584
$data[0,18] = 0.00904779900579266
$data[0,19] = 0.009047799005792664

$data[1,0] = $strs[20]
$data[1,1] = $strs[23]
$data[1,2] = $strs[24]
$data[1,3] = $strs[21]
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 17.586354
$data[1,7] = 52.759062
$data[1,8] = 0.2178245326054132
$data[1,9] = 0.2178245326054132
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.639501
$data[1,13] = 4.918502999999999
$data[1,14] = 0.8947982752789682
$data[1,15] = 0.8947982752789683
$data[1,16] = 28.832844969354
$data[1,17] = 259.495604724186
$data[1,18] = 0.1949090160887711
$data[1,19] = 0.1949090160887711

$data[2,0] = $strs[20]
$data[2,1] = $strs[23]
$data[2,2] = $strs[24]
$data[2,3] = $strs[22]
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 17.586354
$data[2,7] = 52.759062
$data[2,8] = 0.2178245326054132
$data[2,9] = 0.2178245326054132
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.11665
$data[2,13] = 0.34995
$data[2,14] = 0.06366462649994824
$data[2,15] = 0.06366462649994825
$data[2,16] = 2.0514481941
$data[2,17] = 18.4630337469
$data[2,18] = 0.01386771751084943
$data[2,19] = 0.01386771751084943

$data[3,0] = $strs[21]
$data[3,1] = $strs[23]
$data[3,2] = $strs[24]
$data[3,3] = $strs[20]
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 52.27042766666667
$data[3,7] = 156.811283
$data[3,8] = 0.6474213742983183
$data[3,9] = 0.6474213742983183
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.07610666666666667
$data[3,13] = 0.22832
$data[3,14] = 0.04153709822108353
$data[3,15] = 0.04153709822108354
$data[3,16] = 3.978128014951111
$data[3,17] = 35.80315213456
$data[3,18] = 0.02689200521465813
$data[3,19] = 0.02689200521465814

$data[4,0] = $strs[21]
$data[4,1] = $strs[23]
$data[4,2] = $strs[24]
$data[4,3] = $strs[21]
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 52.27042766666667
$data[4,7] = 156.811283
$data[4,8] = 0.6474213742983183
$data[4,9] = 0.6474213742983183
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.639501
$data[4,13] = 4.918502999999999
$data[4,14] = 0.8947982752789682
$data[4,15] = 0.8947982752789683
$data[4,16] = 85.69741842992767
$data[4,17] = 771.2767658693489
$data[4,18] = 0.5793115291008745
$data[4,19] = 0.5793115291008746

$data[5,0] = $strs[21]
$data[5,1] = $strs[23]
$data[5,2] = $strs[24]
$data[5,3] = $strs[22]
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 52.27042766666667
$data[5,7] = 156.811283
$data[5,8] = 0.6474213742983183
$data[5,9] = 0.6474213742983183
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.11665
$data[5,13] = 0.34995
$data[5,14] = 0.06366462649994824
$data[5,15] = 0.06366462649994825
$data[5,16] = 6.097345387316667
$data[5,17] = 54.87610848585
$data[5,18] = 0.04121783998278562
$data[5,19] = 0.04121783998278563

$data[6,0] = $strs[22]
$data[6,1] = $strs[23]
$data[6,2] = $strs[24]
$data[6,3] = $strs[20]
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 10.87955133333333
$data[6,7] = 32.638654
$data[6,8] = 0.1347540930962685
$data[6,9] = 0.1347540930962685
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.07610666666666667
$data[6,13] = 0.22832
$data[6,14] = 0.04153709822108353
$data[6,15] = 0.04153709822108354
$data[6,16] = 0.8280063868088889
$data[6,17] = 7.452057481280002
$data[6,18] = 0.005597294000632738
$data[6,19] = 0.00559729400063274

$data[7,0] = $strs[22]
$data[7,1] = $strs[23]
$data[7,2] = $strs[24]
$data[7,3] = $strs[21]
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 10.87955133333333
$data[7,7] = 32.638654
$data[7,8] = 0.1347540930962685
$data[7,9] = 0.1347540930962685
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.639501
$data[7,13] = 4.918502999999999
$data[7,14] = 0.8947982752789682
$data[7,15] = 0.8947982752789683
$data[7,16] = 17.83703529055133
$data[7,17] = 160.533317614962
$data[7,18] = 0.1205777300893225
$data[7,19] = 0.1205777300893226

$data[8,0] = $strs[22]
$data[8,1] = $strs[23]
$data[8,2] = $strs[24]
$data[8,3] = $strs[22]
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 10.87955133333333
$data[8,7] = 32.638654
$data[8,8] = 0.1347540930962685
$data[8,9] = 0.1347540930962685
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.11665
$data[8,13] = 0.34995
$data[8,14] = 0.06366462649994824
$data[8,15] = 0.06366462649994825
$data[8,16] = 1.269099663033333
$data[8,17] = 11.4218969673
$data[8,18] = 0.008579069006313184
$data[8,19] = 0.008579069006313188

$ws.Range("A2:T10").Value = $data